# Refresh the cryptocurrency price/volume snapshot (Price = column D,
# Volume(1h) = column E) for rows 2-51, as produced by the scheduled
# GitHub Actions scraper run.
#
# Note: several "Price" strings look like plain decimals (e.g. "1.012")
# but must stay TEXT (they're thousand-grouped/display strings, not
# numbers - some even have two dots, e.g. "29.547.90"). A leading
# apostrophe forces Excel to store the literal text instead of
# reinterpreting it as a number; it is a formatting prefix only and is
# not part of the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.547.90'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').Value = '1.920.65'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('D4').Value = '''1.012'
$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').Value = '''326.32'
$ws.Range('E5').Value = '  +0.15%  '

$ws.Range('D6').Value = '''1.011'

$ws.Range('D7').Value = '''0.4817'
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('D8').Value = '''0.4064'
$ws.Range('E8').Value = '  -0.15%  '

$ws.Range('E9').Value = '  +1.17%  '

$ws.Range('D10').Value = '''1.011'
$ws.Range('E10').Value = '  -0.08%  '

$ws.Range('D11').Value = '''23.45'
$ws.Range('E11').Value = '  +0.19%  '

$ws.Range('D12').Value = '1.922.15'
$ws.Range('E12').Value = '  +0.74%  '

$ws.Range('D13').Value = '''6.068'
$ws.Range('E13').Value = '  +0.85%  '

$ws.Range('D14').Value = '''7.254'
$ws.Range('E14').Value = '  +2.10%  '

$ws.Range('D15').Value = '''91.69'
$ws.Range('E15').Value = '  +1.37%  '

$ws.Range('D16').Value = '''0.06877'
$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('E17').Value = '  +0.37%  '

$ws.Range('E18').Value = '  -0.08%  '

$ws.Range('D19').Value = '''17.60'
$ws.Range('E19').Value = '  -0.64%  '

$ws.Range('E20').Value = '  +0.41%  '

$ws.Range('D21').Value = '29.561.43'
$ws.Range('E21').Value = '  +0.00%  '

$ws.Range('D22').Value = '''5.687'

$ws.Range('D23').Value = '''11.89'
$ws.Range('E23').Value = '  +0.34%  '

$ws.Range('D24').Value = '''2.185'
$ws.Range('E24').Value = '  +0.41%  '

$ws.Range('D25').Value = '2.149.53'
$ws.Range('E25').Value = '  +0.57%  '

$ws.Range('D26').Value = '''6.537'
$ws.Range('E26').Value = '  +3.72%  '

$ws.Range('D27').Value = '''155.82'
$ws.Range('E27').Value = '  +0.62%  '

$ws.Range('D28').Value = '''20.05'
$ws.Range('E28').Value = '  -0.04%  '

$ws.Range('D29').Value = '''2.100'
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('D30').Value = '''120.72'
$ws.Range('E30').Value = '  +0.84%  '

$ws.Range('D31').Value = '''1.021'
$ws.Range('E31').Value = '  -0.67%  '

$ws.Range('D32').Value = '''0.09638'
$ws.Range('E32').Value = '  +0.89%  '

$ws.Range('D33').Value = '''5.620'
$ws.Range('E33').Value = '  +1.83%  '

$ws.Range('D34').Value = '''3.556'
$ws.Range('E34').Value = '  -0.20%  '

$ws.Range('D35').Value = '''1.378'
$ws.Range('E35').Value = '  -1.05%  '

$ws.Range('D36').Value = '''0.06360'
$ws.Range('E36').Value = '  +4.16%  '

$ws.Range('D37').Value = '''0.02296'
$ws.Range('E37').Value = '  +1.23%  '

$ws.Range('D38').Value = '''1.194'
$ws.Range('E38').Value = '  +1.37%  '

$ws.Range('D39').Value = '''0.5957'
$ws.Range('E39').Value = '  +0.31%  '

$ws.Range('E40').Value = '  -0.29%  '

$ws.Range('D41').Value = '''7.924'
$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('D42').Value = '''0.1852'
$ws.Range('E42').Value = '  -0.22%  '

$ws.Range('D43').Value = '''2.478'
$ws.Range('E43').Value = '  +0.90%  '

$ws.Range('D44').Value = '''1.285'
$ws.Range('E44').Value = '  +0.14%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('D46').Value = '''0.07495'
$ws.Range('E46').Value = '  -3.11%  '

$ws.Range('D47').Value = '''0.5572'
$ws.Range('E47').Value = '  +0.10%  '

$ws.Range('D48').Value = '''1.949'
$ws.Range('E48').Value = '  +0.07%  '

$ws.Range('D49').Value = '''118.82'
$ws.Range('E49').Value = '  +2.77%  '

$ws.Range('E50').Value = '  +3.53%  '

$ws.Range('D51').Value = '''72.37'
$ws.Range('E51').Value = '  -0.43%  '
